$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.691.58'
$ws.Range('E2').Value = '  +2.14%  '
$ws.Range('D3').Value = '3.004.98'
$ws.Range('E3').Value = '  +1.86%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '512.11'
$ws.Range('E5').Value = '  +5.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.85'
$ws.Range('E6').Value = '  +6.51%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.434'
$ws.Range('E8').Value = '  +4.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.52'
$ws.Range('E9').Value = '  +7.97%  '
$ws.Range('E10').Value = '  +9.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.356'
$ws.Range('E11').Value = '  +2.40%  '
$ws.Range('E12').Value = '  +2.80%  '
$ws.Range('D13').Value = '3.520.06'
$ws.Range('E13').Value = '  +1.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.72'
$ws.Range('E14').Value = '  +5.93%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000158'
$ws.Range('E15').Value = '  +13.80%  '
$ws.Range('D16').Value = '56.748.14'
$ws.Range('E16').Value = '  +2.60%  '
$ws.Range('D17').Value = '3.006.17'
$ws.Range('E17').Value = '  +2.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.96'
$ws.Range('E18').Value = '  +6.98%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.51'
$ws.Range('E19').Value = '  +4.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.85'
$ws.Range('E20').Value = '  +6.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '330.78'
$ws.Range('E21').Value = '  +6.18%  '
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.482'
$ws.Range('E23').Value = '  +5.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.13'
$ws.Range('E24').Value = '  +5.46%  '
$ws.Range('E25').Value = '  +10.40%  '
$ws.Range('D27').Value = '0.0₃0904'
$ws.Range('E27').Value = '  +6.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.73'
$ws.Range('E28').Value = '  +3.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.06'
$ws.Range('E29').Value = '  +9.55%  '
$ws.Range('E30').Value = '  +6.67%  '
$ws.Range('E31').Value = '  +7.74%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.70'
$ws.Range('E32').Value = '  +7.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '153.61'
$ws.Range('E33').Value = '  +4.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.57'
$ws.Range('E34').Value = '  +5.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.67'
$ws.Range('E35').Value = '  +0.75%  '
$ws.Range('E36').Value = '  +0.68%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0678'
$ws.Range('E37').Value = '  +4.91%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '23.75'
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('D39').Value = '3.037.09'
$ws.Range('E39').Value = '  +2.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.94'
$ws.Range('E40').Value = '  +2.08%  '
$ws.Range('E41').Value = '  +0.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.649'
$ws.Range('E42').Value = '  +2.95%  '
$ws.Range('D43').Value = '2.278.23'
$ws.Range('E43').Value = '  +7.85%  '
$ws.Range('E44').Value = '  +0.95%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.68'
$ws.Range('E45').Value = '  +4.78%  '
$ws.Range('E46').Value = '  +3.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.96'
$ws.Range('E47').Value = '  +14.47%  '
$ws.Range('E48').Value = '  +4.95%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.83'
$ws.Range('E49').Value = '  +5.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.43'
$ws.Range('E50').Value = '  +2.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0872'
$ws.Range('E51').Value = '  +6.35%  '
